$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A122").Value = "GO Programmer/Golang Developer"
$ws.Range("B122").Value = "https://www.dice.com/job-detail/596b4b34-7d34-46f1-833c-f5648a5f1157"
$ws.Range("C122").Value = "Mechanicsburg, Pennsylvania"
$ws.Range("D122").Value = "Third Party"
$ws.Range("E122").Value = "USD54 - USD64"
$ws.Range("F122").Value = "Randstad Digital"

$ws.Range("A123").Value = "Golang Developer"
$ws.Range("B123").Value = "https://www.dice.com/job-detail/934b17e2-62d2-4af7-8c8e-6a9c8eab4aab"
$ws.Range("C123").Value = "Remote"
$ws.Range("D123").Value = "Contract"
$ws.Range("E123").Value = "Depends on Experience"
$ws.Range("F123").Value = "SSTech LLC"
